$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add standup data for Sprint 3: new "Tasks Left" values in column C for rows 3 and 4,
# and correct the date in A4.

# Row 3: add C3 = 16, matching formatting of existing B3/C2 cells (style index 3)
$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = 16

# Row 4: update date A4, add C4 = 6 with matching formatting
$ws.Range("A4").Value = 43787
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = 6
